# Applies updated crypto price/volume data to Sheet1, matching the
# latest scrape snapshot (cell values only; no formatting changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.767.01'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.288.96'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '299.70'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.72%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '96.29'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.50%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.511'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -4.13%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '35.47'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.78%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0785'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('E12').Value = '  +0.62%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '17.71'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.97%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.72'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.42%  '
$ws.Range('D15').Value = '2.647.05'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '2.283.41'
$ws.Range('E16').Value = '  +0.84%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.771'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '42.701.70'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.71'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -5.47%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -1.00%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.03'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -2.45%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '67.62'
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '239.52'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.11'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '4.01'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.41%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '24.91'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.12%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '165.88'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -1.96%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '32.67'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('E35').Value = '  -4.44%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '16.85'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -8.09%  '
$ws.Range('E37').Value = '  -1.69%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0682'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('E40').Value = '  -3.50%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.109'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.32%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.69'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').Value = '2.009.66'
$ws.Range('E43').Value = '  +0.83%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0280'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.90%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '10.04'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.63%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.12'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.26%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '17.07'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.99%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.76'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.77%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.93'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.70%  '
$ws.Range('D50').Value = '2.515.15'
$ws.Range('E50').Value = '  -0.83%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '52.88'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.19%  '
